# Update the cryptos list with the latest scraped Coinranking data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "Price" cell (column D) as literal text even when the
# string looks numeric (e.g. "1.005"), matching the inline-string source
# data -- a leading apostrophe forces text entry, then the style is reset
# to Normal so no stray number-format style sticks to the cell.
function Set-PriceText($ref, $text) {
    $c = $ws.Range($ref)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}


# Row 2
Set-PriceText "D2" '30.272.32'
$ws.Range("E2").Value = '  +1.36%  '

# Row 3
Set-PriceText "D3" '2.087.94'
$ws.Range("E3").Value = '  -1.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
Set-PriceText "D5" '341.20'
$ws.Range("E5").Value = '  -1.97%  '

# Row 6
Set-PriceText "D6" '1.005'
$ws.Range("E6").Value = '  -0.18%  '

# Row 7
$ws.Range("E7").Value = '  +1.92%  '

# Row 8
Set-PriceText "D8" '0.4390'
$ws.Range("E8").Value = '  -1.81%  '

# Row 9
Set-PriceText "D9" '54.91'
$ws.Range("E9").Value = '  +1.50%  '

# Row 10
Set-PriceText "D10" '0.09353'
$ws.Range("E10").Value = '  -0.06%  '

# Row 11
Set-PriceText "D11" '1.176'
$ws.Range("E11").Value = '  -0.59%  '

# Row 12
Set-PriceText "D12" '24.50'
$ws.Range("E12").Value = '  -2.77%  '

# Row 13
Set-PriceText "D13" '8.499'
$ws.Range("E13").Value = '  +1.30%  '

# Row 14
Set-PriceText "D14" '6.866'
$ws.Range("E14").Value = '  +0.36%  '

# Row 15
Set-PriceText "D15" '2.037.02'
$ws.Range("E15").Value = '  -2.85%  '

# Row 16
Set-PriceText "D16" '101.62'
$ws.Range("E16").Value = '  -1.03%  '

# Row 17
$ws.Range("E17").Value = '  -0.93%  '

# Row 18
Set-PriceText "D18" '1.006'

# Row 19
Set-PriceText "D19" '21.03'
$ws.Range("E19").Value = '  -2.57%  '

# Row 20
Set-PriceText "D20" '0.06712'
$ws.Range("E20").Value = '  +0.57%  '

# Row 21
Set-PriceText "D21" '6.293'
$ws.Range("E21").Value = '  -0.20%  '

# Row 22
Set-PriceText "D22" '1.004'
$ws.Range("E22").Value = '  -0.22%  '

# Row 23
Set-PriceText "D23" '30.281.19'
$ws.Range("E23").Value = '  +1.27%  '

# Row 24
Set-PriceText "D24" '12.40'
$ws.Range("E24").Value = '  -2.47%  '

# Row 25
Set-PriceText "D25" '2.325'
$ws.Range("E25").Value = '  -0.20%  '

# Row 26
$ws.Range("E26").Value = '  -1.71%  '

# Row 27
Set-PriceText "D27" '6.826'
$ws.Range("E27").Value = '  +6.33%  '

# Row 28
Set-PriceText "D28" '162.65'
$ws.Range("E28").Value = '  -0.01%  '

# Row 29
Set-PriceText "D29" '2.482'
$ws.Range("E29").Value = '  -2.78%  '

# Row 30
Set-PriceText "D30" '133.71'
$ws.Range("E30").Value = '  -0.32%  '

# Row 31
Set-PriceText "D31" '1.128'
$ws.Range("E31").Value = '  -2.52%  '

# Row 32
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-PriceText "D32" '1.663'
$ws.Range("E32").Value = '  -7.40%  '

# Row 33
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-PriceText "D33" '0.1048'
$ws.Range("E33").Value = '  -0.65%  '

# Row 34
Set-PriceText "D34" '6.258'
$ws.Range("E34").Value = '  +0.13%  '

# Row 35
Set-PriceText "D35" '3.915'
$ws.Range("E35").Value = '  -1.44%  '

# Row 36
Set-PriceText "D36" '0.02609'
$ws.Range("E36").Value = '  +0.46%  '

# Row 37
Set-PriceText "D37" '9.903'
$ws.Range("E37").Value = '  -9.17%  '

# Row 38
Set-PriceText "D38" '0.06758'
$ws.Range("E38").Value = '  -0.63%  '

# Row 39
Set-PriceText "D39" '12.57'
$ws.Range("E39").Value = '  -1.10%  '

# Row 40
Set-PriceText "D40" '1.343'
$ws.Range("E40").Value = '  -0.38%  '

# Row 41
Set-PriceText "D41" '0.6927'
$ws.Range("E41").Value = '  -1.38%  '

# Row 42
Set-PriceText "D42" '0.2205'
$ws.Range("E42").Value = '  -1.84%  '

# Row 43
Set-PriceText "D43" '0.6743'
$ws.Range("E43").Value = '  -1.61%  '

# Row 44
Set-PriceText "D44" '2.384'
$ws.Range("E44").Value = '  +0.85%  '

# Row 45
Set-PriceText "D45" '14.29'
$ws.Range("E45").Value = '  -1.34%  '

# Row 46
$ws.Range("E46").Value = '  -0.22%  '

# Row 47
Set-PriceText "D47" '1.290'
$ws.Range("E47").Value = '  +5.74%  '

# Row 48
Set-PriceText "D48" '3.639'
$ws.Range("E48").Value = '  +0.05%  '

# Row 49
Set-PriceText "D49" '0.00000000345'
$ws.Range("E49").Value = '  -3.66%  '

# Row 50
Set-PriceText "D50" '1.210'
$ws.Range("E50").Value = '  +2.33%  '

# Row 51
Set-PriceText "D51" '1.209'
$ws.Range("E51").Value = '  -1.39%  '
